$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value2 = [double]"3"
$ws.Range("G2").Value2 = [double]"85.28734300000001"
$ws.Range("H2").Value2 = [double]"255.862029"
$ws.Range("I2").Value2 = [double]"0.7886846474407839"
$ws.Range("J2").Value2 = [double]"0.7886846474407838"
$ws.Range("K2").Value2 = [double]"3"
$ws.Range("M2").Value2 = [double]"85.28734300000001"
$ws.Range("N2").Value2 = [double]"255.862029"
$ws.Range("O2").Value2 = [double]"0.7886846474407839"
$ws.Range("P2").Value2 = [double]"0.7886846474407838"
$ws.Range("Q2").Value2 = [double]"7273.93087599965"
$ws.Range("R2").Value2 = [double]"65465.37788399684"
$ws.Range("S2").Value2 = [double]"0.6220234731087937"
$ws.Range("T2").Value2 = [double]"0.6220234731087935"

# Row 3
$ws.Range("E3").Value2 = [double]"3"
$ws.Range("G3").Value2 = [double]"85.28734300000001"
$ws.Range("H3").Value2 = [double]"255.862029"
$ws.Range("I3").Value2 = [double]"0.7886846474407839"
$ws.Range("J3").Value2 = [double]"0.7886846474407838"
$ws.Range("K3").Value2 = [double]"3"
$ws.Range("M3").Value2 = [double]"19.11746166666667"
$ws.Range("N3").Value2 = [double]"57.352385"
$ws.Range("O3").Value2 = [double]"0.1767864724609571"
$ws.Range("P3").Value2 = [double]"0.1767864724609571"
$ws.Range("Q3").Value2 = [double]"1630.477510454352"
$ws.Range("R3").Value2 = [double]"14674.29759408916"
$ws.Range("S3").Value2 = [double]"0.1394287767051698"
$ws.Range("T3").Value2 = [double]"0.1394287767051698"

# Row 4
$ws.Range("E4").Value2 = [double]"3"
$ws.Range("G4").Value2 = [double]"85.28734300000001"
$ws.Range("H4").Value2 = [double]"255.862029"
$ws.Range("I4").Value2 = [double]"0.7886846474407839"
$ws.Range("J4").Value2 = [double]"0.7886846474407838"
$ws.Range("K4").Value2 = [double]"3"
$ws.Range("M4").Value2 = [double]"0.5740886666666668"
$ws.Range("N4").Value2 = [double]"1.722266"
$ws.Range("O4").Value2 = [double]"0.005308817249351405"
$ws.Range("P4").Value2 = [double]"0.005308817249351405"
$ws.Range("Q4").Value2 = [double]"48.96249702641268"
$ws.Range("R4").Value2 = [double]"440.6624732377141"
$ws.Range("S4").Value2 = [double]"0.004186982660632265"
$ws.Range("T4").Value2 = [double]"0.004186982660632264"

# Row 5
$ws.Range("E5").Value2 = [double]"3"
$ws.Range("G5").Value2 = [double]"85.28734300000001"
$ws.Range("H5").Value2 = [double]"255.862029"
$ws.Range("I5").Value2 = [double]"0.7886846474407839"
$ws.Range("J5").Value2 = [double]"0.7886846474407838"
$ws.Range("K5").Value2 = [double]"3"
$ws.Range("M5").Value2 = [double]"3.15982"
$ws.Range("N5").Value2 = [double]"9.479460000000001"
$ws.Range("O5").Value2 = [double]"0.02922006284890758"
$ws.Range("P5").Value2 = [double]"0.02922006284890758"
$ws.Range("Q5").Value2 = [double]"269.4926521582601"
$ws.Range("R5").Value2 = [double]"2425.433869424341"
$ws.Range("S5").Value2 = [double]"0.02304541496618822"
$ws.Range("T5").Value2 = [double]"0.02304541496618822"

# Row 6
$ws.Range("E6").Value2 = [double]"3"
$ws.Range("G6").Value2 = [double]"19.11746166666667"
$ws.Range("H6").Value2 = [double]"57.352385"
$ws.Range("I6").Value2 = [double]"0.1767864724609571"
$ws.Range("J6").Value2 = [double]"0.1767864724609571"
$ws.Range("K6").Value2 = [double]"3"
$ws.Range("M6").Value2 = [double]"85.28734300000001"
$ws.Range("N6").Value2 = [double]"255.862029"
$ws.Range("O6").Value2 = [double]"0.7886846474407839"
$ws.Range("P6").Value2 = [double]"0.7886846474407838"
$ws.Range("Q6").Value2 = [double]"1630.477510454352"
$ws.Range("R6").Value2 = [double]"14674.29759408916"
$ws.Range("S6").Value2 = [double]"0.1394287767051698"
$ws.Range("T6").Value2 = [double]"0.1394287767051698"

# Row 7
$ws.Range("E7").Value2 = [double]"3"
$ws.Range("G7").Value2 = [double]"19.11746166666667"
$ws.Range("H7").Value2 = [double]"57.352385"
$ws.Range("I7").Value2 = [double]"0.1767864724609571"
$ws.Range("J7").Value2 = [double]"0.1767864724609571"
$ws.Range("K7").Value2 = [double]"3"
$ws.Range("M7").Value2 = [double]"19.11746166666667"
$ws.Range("N7").Value2 = [double]"57.352385"
$ws.Range("O7").Value2 = [double]"0.1767864724609571"
$ws.Range("P7").Value2 = [double]"0.1767864724609571"
$ws.Range("Q7").Value2 = [double]"365.4773405764695"
$ws.Range("R7").Value2 = [double]"3289.296065188225"
$ws.Range("S7").Value2 = [double]"0.03125345684518874"
$ws.Range("T7").Value2 = [double]"0.03125345684518872"

# Row 8
$ws.Range("E8").Value2 = [double]"3"
$ws.Range("G8").Value2 = [double]"19.11746166666667"
$ws.Range("H8").Value2 = [double]"57.352385"
$ws.Range("I8").Value2 = [double]"0.1767864724609571"
$ws.Range("J8").Value2 = [double]"0.1767864724609571"
$ws.Range("K8").Value2 = [double]"3"
$ws.Range("M8").Value2 = [double]"0.5740886666666668"
$ws.Range("N8").Value2 = [double]"1.722266"
$ws.Range("O8").Value2 = [double]"0.005308817249351405"
$ws.Range("P8").Value2 = [double]"0.005308817249351405"
$ws.Range("Q8").Value2 = [double]"10.97511807826778"
$ws.Range("R8").Value2 = [double]"98.77606270441001"
$ws.Range("S8").Value2 = [double]"0.000938527074452716"
$ws.Range("T8").Value2 = [double]"0.0009385270744527158"

# Row 9
$ws.Range("E9").Value2 = [double]"3"
$ws.Range("G9").Value2 = [double]"19.11746166666667"
$ws.Range("H9").Value2 = [double]"57.352385"
$ws.Range("I9").Value2 = [double]"0.1767864724609571"
$ws.Range("J9").Value2 = [double]"0.1767864724609571"
$ws.Range("K9").Value2 = [double]"3"
$ws.Range("M9").Value2 = [double]"3.15982"
$ws.Range("N9").Value2 = [double]"9.479460000000001"
$ws.Range("O9").Value2 = [double]"0.02922006284890758"
$ws.Range("P9").Value2 = [double]"0.02922006284890758"
$ws.Range("Q9").Value2 = [double]"60.40773772356668"
$ws.Range("R9").Value2 = [double]"543.6696395121"
$ws.Range("S9").Value2 = [double]"0.005165711836145835"
$ws.Range("T9").Value2 = [double]"0.005165711836145834"

# Row 10
$ws.Range("E10").Value2 = [double]"3"
$ws.Range("G10").Value2 = [double]"0.5740886666666668"
$ws.Range("H10").Value2 = [double]"1.722266"
$ws.Range("I10").Value2 = [double]"0.005308817249351405"
$ws.Range("J10").Value2 = [double]"0.005308817249351405"
$ws.Range("K10").Value2 = [double]"3"
$ws.Range("M10").Value2 = [double]"85.28734300000001"
$ws.Range("N10").Value2 = [double]"255.862029"
$ws.Range("O10").Value2 = [double]"0.7886846474407839"
$ws.Range("P10").Value2 = [double]"0.7886846474407838"
$ws.Range("Q10").Value2 = [double]"48.96249702641268"
$ws.Range("R10").Value2 = [double]"440.6624732377141"
$ws.Range("S10").Value2 = [double]"0.004186982660632265"
$ws.Range("T10").Value2 = [double]"0.004186982660632264"

# Row 11
$ws.Range("E11").Value2 = [double]"3"
$ws.Range("G11").Value2 = [double]"0.5740886666666668"
$ws.Range("H11").Value2 = [double]"1.722266"
$ws.Range("I11").Value2 = [double]"0.005308817249351405"
$ws.Range("J11").Value2 = [double]"0.005308817249351405"
$ws.Range("K11").Value2 = [double]"3"
$ws.Range("M11").Value2 = [double]"19.11746166666667"
$ws.Range("N11").Value2 = [double]"57.352385"
$ws.Range("O11").Value2 = [double]"0.1767864724609571"
$ws.Range("P11").Value2 = [double]"0.1767864724609571"
$ws.Range("Q11").Value2 = [double]"10.97511807826778"
$ws.Range("R11").Value2 = [double]"98.77606270441001"
$ws.Range("S11").Value2 = [double]"0.000938527074452716"
$ws.Range("T11").Value2 = [double]"0.0009385270744527158"

# Row 12
$ws.Range("E12").Value2 = [double]"3"
$ws.Range("G12").Value2 = [double]"0.5740886666666668"
$ws.Range("H12").Value2 = [double]"1.722266"
$ws.Range("I12").Value2 = [double]"0.005308817249351405"
$ws.Range("J12").Value2 = [double]"0.005308817249351405"
$ws.Range("K12").Value2 = [double]"3"
$ws.Range("M12").Value2 = [double]"0.5740886666666668"
$ws.Range("N12").Value2 = [double]"1.722266"
$ws.Range("O12").Value2 = [double]"0.005308817249351405"
$ws.Range("P12").Value2 = [double]"0.005308817249351405"
$ws.Range("Q12").Value2 = [double]"0.3295777971951113"
$ws.Range("R12").Value2 = [double]"2.966200174756001"
$ws.Range("S12").Value2 = [double]"2.818354058701101E-05"
$ws.Range("T12").Value2 = [double]"2.818354058701101E-05"

# Row 13
$ws.Range("E13").Value2 = [double]"3"
$ws.Range("G13").Value2 = [double]"0.5740886666666668"
$ws.Range("H13").Value2 = [double]"1.722266"
$ws.Range("I13").Value2 = [double]"0.005308817249351405"
$ws.Range("J13").Value2 = [double]"0.005308817249351405"
$ws.Range("K13").Value2 = [double]"3"
$ws.Range("M13").Value2 = [double]"3.15982"
$ws.Range("N13").Value2 = [double]"9.479460000000001"
$ws.Range("O13").Value2 = [double]"0.02922006284890758"
$ws.Range("P13").Value2 = [double]"0.02922006284890758"
$ws.Range("Q13").Value2 = [double]"1.814016850706667"
$ws.Range("R13").Value2 = [double]"16.32615165636"
$ws.Range("S13").Value2 = [double]"0.0001551239736794127"
$ws.Range("T13").Value2 = [double]"0.0001551239736794127"

# Row 14
$ws.Range("E14").Value2 = [double]"3"
$ws.Range("G14").Value2 = [double]"3.15982"
$ws.Range("H14").Value2 = [double]"9.479460000000001"
$ws.Range("I14").Value2 = [double]"0.02922006284890758"
$ws.Range("J14").Value2 = [double]"0.02922006284890758"
$ws.Range("K14").Value2 = [double]"3"
$ws.Range("M14").Value2 = [double]"85.28734300000001"
$ws.Range("N14").Value2 = [double]"255.862029"
$ws.Range("O14").Value2 = [double]"0.7886846474407839"
$ws.Range("P14").Value2 = [double]"0.7886846474407838"
$ws.Range("Q14").Value2 = [double]"269.4926521582601"
$ws.Range("R14").Value2 = [double]"2425.433869424341"
$ws.Range("S14").Value2 = [double]"0.02304541496618822"
$ws.Range("T14").Value2 = [double]"0.02304541496618822"

# Row 15
$ws.Range("E15").Value2 = [double]"3"
$ws.Range("G15").Value2 = [double]"3.15982"
$ws.Range("H15").Value2 = [double]"9.479460000000001"
$ws.Range("I15").Value2 = [double]"0.02922006284890758"
$ws.Range("J15").Value2 = [double]"0.02922006284890758"
$ws.Range("K15").Value2 = [double]"3"
$ws.Range("M15").Value2 = [double]"19.11746166666667"
$ws.Range("N15").Value2 = [double]"57.352385"
$ws.Range("O15").Value2 = [double]"0.1767864724609571"
$ws.Range("P15").Value2 = [double]"0.1767864724609571"
$ws.Range("Q15").Value2 = [double]"60.40773772356668"
$ws.Range("R15").Value2 = [double]"543.6696395121"
$ws.Range("S15").Value2 = [double]"0.005165711836145835"
$ws.Range("T15").Value2 = [double]"0.005165711836145834"

# Row 16
$ws.Range("E16").Value2 = [double]"3"
$ws.Range("G16").Value2 = [double]"3.15982"
$ws.Range("H16").Value2 = [double]"9.479460000000001"
$ws.Range("I16").Value2 = [double]"0.02922006284890758"
$ws.Range("J16").Value2 = [double]"0.02922006284890758"
$ws.Range("K16").Value2 = [double]"3"
$ws.Range("M16").Value2 = [double]"0.5740886666666668"
$ws.Range("N16").Value2 = [double]"1.722266"
$ws.Range("O16").Value2 = [double]"0.005308817249351405"
$ws.Range("P16").Value2 = [double]"0.005308817249351405"
$ws.Range("Q16").Value2 = [double]"1.814016850706667"
$ws.Range("R16").Value2 = [double]"16.32615165636"
$ws.Range("S16").Value2 = [double]"0.0001551239736794127"
$ws.Range("T16").Value2 = [double]"0.0001551239736794127"

# Row 17
$ws.Range("E17").Value2 = [double]"3"
$ws.Range("G17").Value2 = [double]"3.15982"
$ws.Range("H17").Value2 = [double]"9.479460000000001"
$ws.Range("I17").Value2 = [double]"0.02922006284890758"
$ws.Range("J17").Value2 = [double]"0.02922006284890758"
$ws.Range("K17").Value2 = [double]"3"
$ws.Range("M17").Value2 = [double]"3.15982"
$ws.Range("N17").Value2 = [double]"9.479460000000001"
$ws.Range("O17").Value2 = [double]"0.02922006284890758"
$ws.Range("P17").Value2 = [double]"0.02922006284890758"
$ws.Range("Q17").Value2 = [double]"9.984462432400003"
$ws.Range("R17").Value2 = [double]"89.86016189160003"
$ws.Range("S17").Value2 = [double]"0.0008538120728941088"
$ws.Range("T17").Value2 = [double]"0.0008538120728941088"
